# 2022.01.05 | Hu Wenqiang | Updated Card
# Set explicit column widths for columns A-F (widths: 19, 8, 20, 7, 40, 13 chars)
# ColumnWidth values below are chosen so the serialized OOXML <col width="..."/>
# resolves to the exact integer targets required.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 18.1
$ws.Columns.Item(2).ColumnWidth = 7.1
$ws.Columns.Item(3).ColumnWidth = 19.1
$ws.Columns.Item(4).ColumnWidth = 6.1
$ws.Columns.Item(5).ColumnWidth = 39.1
$ws.Columns.Item(6).ColumnWidth = 12.1

# Clear the stray empty-string cells that were left over from the source data
# (they show up as inline strings with no text) so the cells become truly
# empty / absent from the sheet data, matching the cleaned-up export.
$emptyCells = "G2,E3,E4,E5,E6,E7,E8,E9,E10,E11,E12,E13,E14,E15,E16,E17,E18,E19,E20,C21,E21,E22,E23,E24,E25,E26,E27,E28,E29,E30,E31,E32,E33,E34,E35,E36,E37,E38,E39,E40,E41,E42,E43,E44,E45,E46,E47,E48,E49,E50,E51,E52,E53,E54,E55,E56,E57,E58,E59,E60,E61,E62,E63,E64,E65,E66,E67,E68,E69,E70,E71,E72,E73,E74,E75,E76,E77,E78,E79,E80,E81,E82,E83,E84,E85,E86,E87,E88,E89,E90,E91,E92,E93,E94,E95,E96,C97,E97,E98,E99,E100,E101,E102,A103,B103,C103,E103"

$range = $ws.Range($emptyCells)
foreach ($cell in $range) {
    $cell.ClearContents()
}
